$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "datetimeFigureOut" field text: 2022/1/23 -> 2022/1/24
#    Present once on the slide master and once on each of the 11 layouts.
# ---------------------------------------------------------------------------
function Update-DateText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2022/1/23") {
                $tr.Text = "2022/1/24"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateText $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateText $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2) Fig. 14 diagram (slides 2 & 3): shift every shape 1678375 EMU to the
#    right so the inserted explanatory video fits beside it without overlap.
#    (Y position is unchanged.) Target X positions are keyed by Shape.Id so
#    the mapping is robust to any reordering of the shape collection.
# ---------------------------------------------------------------------------
$EMU_PER_POINT = 12700.0

$slide2TargetX = @{
  "121" = 3355202
  "55" = 4522361
  "54" = 4522360
  "5" = 3308107
  "7" = 4234288
  "20" = 5664442
  "35" = 4830346
  "38" = 5378921
  "39" = 5093401
  "40" = 4807881
  "41" = 4522360
  "43" = 4234288
  "44" = 4234288
  "45" = 4234288
  "47" = 4234288
  "48" = 4234288
  "56" = 4522153
  "57" = 4807663
  "58" = 4808616
  "59" = 5092732
  "60" = 5092455
  "61" = 5376891
  "53" = 4608867
  "62" = 6340279
  "63" = 6670800
  "64" = 5046382
  "65" = 5046382
  "66" = 5046382
  "67" = 4759947
  "68" = 4759947
  "69" = 4759947
  "70" = 4477229
  "71" = 4477229
  "72" = 4477229
  "73" = 5333677
  "74" = 5333677
  "75" = 5333677
  "76" = 4759947
  "77" = 4759947
  "78" = 4477229
  "79" = 5046382
  "80" = 5617289
  "81" = 5617289
  "82" = 5617289
  "83" = 5617289
  "84" = 5617289
  "85" = 5617289
  "86" = 4477229
  "87" = 4759947
  "88" = 4477229
  "89" = 5333677
  "90" = 5333677
  "91" = 5333677
  "93" = 5046382
  "94" = 5046382
  "96" = 6433930
  "97" = 6670800
  "98" = 6433930
  "99" = 6670800
  "100" = 6670800
  "101" = 6670800
  "102" = 6433930
  "103" = 6439570
  "104" = 3355202
  "106" = 5333677
  "2" = 5351968
  "92" = 5064420
  "95" = 4776872
  "105" = 5625375
  "107" = 5318304
  "108" = 4776872
  "109" = 4492507
  "110" = 4210267
  "111" = 4492507
  "112" = 5022742
  "113" = 5318304
  "114" = 4210267
  "115" = 4210267
  "116" = 4210267
  "117" = 4735194
  "118" = 4735194
  "119" = 4735194
  "120" = 5022742
}

$slide3TargetX = @{
  "125" = 3355202
  "55" = 4522361
  "54" = 4522360
  "5" = 3308107
  "7" = 4234288
  "20" = 5664442
  "35" = 4830346
  "38" = 5378921
  "39" = 5093401
  "40" = 4807881
  "41" = 4522360
  "43" = 4234288
  "44" = 4234288
  "45" = 4234288
  "47" = 4234288
  "48" = 4234288
  "56" = 4522153
  "57" = 4807663
  "58" = 4808616
  "59" = 5092732
  "60" = 5092455
  "61" = 5376891
  "53" = 4608867
  "62" = 6340279
  "63" = 6670800
  "64" = 5046382
  "65" = 5046382
  "66" = 5046382
  "67" = 4759947
  "68" = 4759947
  "69" = 4759947
  "70" = 4477229
  "71" = 4477229
  "72" = 4477229
  "73" = 5333677
  "74" = 5333677
  "75" = 5333677
  "76" = 4759947
  "77" = 4759947
  "78" = 4477229
  "79" = 5046382
  "80" = 5617289
  "81" = 5617289
  "82" = 5617289
  "83" = 5617289
  "84" = 5617289
  "85" = 5617289
  "86" = 4477229
  "87" = 4759947
  "88" = 4477229
  "89" = 5333677
  "90" = 5333677
  "91" = 5333677
  "93" = 5046382
  "94" = 5046382
  "96" = 6433930
  "97" = 6670800
  "98" = 6433930
  "99" = 6670800
  "100" = 6670800
  "101" = 6670800
  "102" = 6433930
  "103" = 6439570
  "104" = 3355202
  "106" = 5333677
  "92" = 5617289
  "95" = 5351968
  "105" = 5064420
  "107" = 4776872
  "108" = 5625375
  "109" = 5318304
  "110" = 4776872
  "111" = 4492507
  "112" = 4210267
  "113" = 4492507
  "114" = 5022742
  "115" = 5318304
  "116" = 4210267
  "117" = 4210267
  "118" = 4210267
  "119" = 4735194
  "120" = 4735194
  "121" = 4735194
  "122" = 5022742
  "3" = 5402254
  "123" = 5588390
  "124" = 5625375
}

function Move-Shapes($slide, $targetMap) {
    $shapes = $slide.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $key = [string]$sh.Id
        if ($targetMap.ContainsKey($key)) {
            $sh.Left = $targetMap[$key] / $EMU_PER_POINT
        }
    }
}

Move-Shapes $p.Slides.Item(2) $slide2TargetX
Move-Shapes $p.Slides.Item(3) $slide3TargetX

Write-Host "done"
